$wb = $excel.ActiveWorkbook

# --- TestCase sheet: insert a new row 14 with "Search_For_Instruction" ---
$ws = $wb.Worksheets.Item("TestCase")
$ws.Activate()

$ws.Rows.Item(14).Insert()
$ws.Cells.Item(14, 2).Value = "Search_For_Instruction"

# Highlight B11 and B12 the same way the rows above them (B2:B10) are
# highlighted, by copying the format from B2.
$ws.Cells.Item(2, 2).Copy()
$ws.Range("B11:B12").PasteSpecial(-4122)

# Move the selection on the TestCase sheet to B4.
$ws.Range("B4").Select()

# --- TestSchedule sheet: move the selection to B2 ---
$wsSchedule = $wb.Worksheets.Item("TestSchedule")
$wsSchedule.Activate()
$wsSchedule.Range("B2").Select()

# Restore TestCase as the active sheet/tab (matches original active tab).
$ws.Activate()
